$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ E=3; G=51.56869633333334; H=154.706089; I=0.1855839901450455; J=0.1855839901450455; K=3; M=32.658611; N=97.97583299999999; O=0.8917945513604759; P=0.8917945513604759; Q=1684.16199332746; R=15157.45793994714; S=0.1655027912310878; T=0.1655027912310878 }
    3  = @{ E=3; G=51.56869633333334; H=154.706089; I=0.1855839901450455; J=0.1855839901450455; K=3; M=0.3385316666666667; N=1.015595; O=0.009244137657793045; P=0.009244137657793047; Q=17.45763671755056; R=157.118730457955; S=0.001715563951983308; T=0.001715563951983308 }
    4  = @{ E=3; G=51.56869633333334; H=154.706089; I=0.1855839901450455; J=0.1855839901450455; K=3; M=3.624084666666667; N=10.872254; O=0.09896131098173097; P=0.09896131098173098; Q=186.8893216616229; R=1682.003894954606; S=0.01836563496197434; T=0.01836563496197434 }
    5  = @{ E=3; G=165.03522; H=495.1056600000001; I=0.5939241598059933; J=0.5939241598059933; K=3; M=32.658611; N=97.97583299999999; O=0.8917945513604759; P=0.8917945513604759; Q=5389.82105127942; R=48508.38946151478; S=0.5296583296363334; T=0.5296583296363334 }
    6  = @{ E=3; G=165.03522; H=495.1056600000001; I=0.5939241598059933; J=0.5939241598059933; K=3; M=0.3385316666666667; N=1.015595; O=0.009244137657793045; P=0.009244137657793047; Q=55.86964808530001; R=502.8268327677001; S=0.005490316691535677; T=0.005490316691535679 }
    7  = @{ E=3; G=165.03522; H=495.1056600000001; I=0.5939241598059933; J=0.5939241598059933; K=3; M=3.624084666666667; N=10.872254; O=0.09896131098173097; P=0.09896131098173098; Q=598.10161026196; R=5382.91449235764; S=0.05877551347812419; T=0.05877551347812419 }
    8  = @{ E=3; G=61.26863233333334; H=183.805897; I=0.2204918500489612; J=0.2204918500489612; K=3; M=32.658611; N=97.97583299999999; O=0.8917945513604759; P=0.8917945513604759; Q=2000.948429876356; R=18008.5358688872; S=0.1966334304930547; T=0.1966334304930547 }
    9  = @{ E=3; G=61.26863233333334; H=183.805897; I=0.2204918500489612; J=0.2204918500489612; K=3; M=0.3385316666666667; N=1.015595; O=0.009244137657793045; P=0.009244137657793047; Q=20.74137221819056; R=186.672349963715; S=0.002038257014274059; T=0.00203825701427406 }
    10 = @{ E=3; G=61.26863233333334; H=183.805897; I=0.2204918500489612; J=0.2204918500489612; K=3; M=3.624084666666667; N=10.872254; O=0.09896131098173097; P=0.09896131098173098; Q=222.0427109868709; R=1998.384398881838; S=0.02182016254163244; T=0.02182016254163244 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
